$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 239, pushing the existing rows 239..341 down to 240..342.
$ws.Rows(239).Insert()

# Populate the newly inserted row 239 with the new data record.
$ws.Range("A239").Value = 10
$ws.Range("B239").Value = "Vega Modelo de Temuco"
$ws.Range("C239").Value = "La Araucanía"
$ws.Range("D239").Value = 44795
$ws.Range("E239").Value = 9
$ws.Range("F239").Value = 100114013
$ws.Range("G239").Value = "Zanahoria"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 80
$ws.Range("K239").Value = 12000
$ws.Range("L239").Value = 12000
$ws.Range("M239").Value = 12000
$ws.Range("N239").Value = "`$/saco 25 kilos"
$ws.Range("O239").Value = "Región de La Araucanía"
$ws.Range("P239").Value = 480
$ws.Range("Q239").Value = 25
$ws.Range("R239").Value = "Hortaliza"
